# This workbook contains ligand/receptor pair statistics per sending cluster.
# The data for the "ECs" sending cluster (row 2) was removed entirely, and the
# specificity figures for the remaining "FAPs" / "MuSCs" rows were recomputed
# against the smaller (2-row) dataset (values supplied by the new TPM run).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ECs" row (old row 2). This shifts the old row3 (FAPs) and
# old row4 (MuSCs) up to become the new row2 and row3, and also drops the
# now-unused "ECs" shared string automatically.
$ws.Rows(2).Delete()

# Update the recalculated specificity columns for the new row 2 (FAPs -> Il18/Il1rapl1 -> MuSCs)
$ws.Range("I2").Value2 = 0.8792824445743517
$ws.Range("J2").Value2 = 0.8792824445743516
$ws.Range("S2").Value2 = 0.8792824445743517
$ws.Range("T2").Value2 = 0.8792824445743516

# Update the recalculated specificity columns for the new row 3 (MuSCs -> Il18/Il1rapl1 -> MuSCs)
$ws.Range("I3").Value2 = 0.1207175554256484
$ws.Range("J3").Value2 = 0.1207175554256484
$ws.Range("S3").Value2 = 0.1207175554256484
$ws.Range("T3").Value2 = 0.1207175554256484
